# Fill in the "answer" column (C) of worksheet "3" (the third exam-attempt
# sheet) with the submitted answers. Column D already holds a shared
# IF(B=C,TRUE,FALSE) formula that will recalculate automatically once C is
# populated, which also updates the summary COUNTIF/percentage cells
# further down the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3")

$ws.Range("C2").Value = 'Use the 802.1x protocol'
$ws.Range("C3").Value = 'Place a front-end web server in a demilitarized zone that only handles external web traffic'
$ws.Range("C4").Value = 'ESP '
$ws.Range("C5").Value = 'Source IP address '
$ws.Range("C6").Value = 'Encrypting the JTAG interface'
$ws.Range("C7").Value = 'Test password strength, brute-force encrypted or hashed passwords, and crack passwords via dictionary attacks'
$ws.Range("C8").Value = 'Grey-box'
$ws.Range("C9").Value = 'BBProxy'
$ws.Range("C10").Value = 'Stretching'
$ws.Range("C11").Value = 'NIDS'
$ws.Range("C12").Value = 'Administrative safeguards '
$ws.Range("C13").Value = 'Private'
$ws.Range("C14").Value = 'Default Credential'
$ws.Range("C15").Value = 'Nothing because the password file does not contain the passwords themselves'
$ws.Range("C16").Value = 2
$ws.Range("C17").Value = 'Watering Hole'
$ws.Range("C18").Value = 'Rules of Engagement'
$ws.Range("C19").Value = 'IP Spoofing'
$ws.Range("C20").Value = 'Display of the contents of the passwd file'
$ws.Range("C21").Value = 'Spectre and Meltdown'
$ws.Range("C22").Value = 'C'
$ws.Range("C23").Value = 'Periodic security checks and audits are required. Access to network elements should be provided by user IDs with strong passwords'
$ws.Range("C24").Value = 'Open-source intelligence'
$ws.Range("C25").Value = 'DNS spoofing '
$ws.Range("C27").Value = 'Nslookup'
$ws.Range("C28").Value = 'Sniffers operate on Layer 2 of the OSI model'
$ws.Range("C29").Value = 'The unauthorized disclosure, modification, or use of sensitive data'
$ws.Range("C30").Value = 'host -t a resolveddomain.com'
$ws.Range("C31").Value = 'Windows'
$ws.Range("C32").Value = 'Non-repudiation'
$ws.Range("C33").Value = 'Open mail relay'
$ws.Range("C34").Value = 'Port scanning, banner grabbing, service identification'
$ws.Range("C36").Value = 'http-methods'
$ws.Range("C37").Value = 'Fuzz Testing'
$ws.Range("C38").Value = 'SYN'
$ws.Range("C39").Value = 'Determine the impact of enabling the audit feature'
$ws.Range("C40").Value = 'An authentication system that creates one-time passwords that are encrypted with secret keys'
$ws.Range("C41").Value = 'Brute force'
$ws.Range("C42").Value = 'Metasploit'
$ws.Range("C43").Value = 'Network layer '
$ws.Range("C44").Value = 'Crypter'
$ws.Range("C45").Value = 'Netbios enumeration'
$ws.Range("C46").Value = 'Markov Chain'
$ws.Range("C47").Value = 'Protocol analyzer'
$ws.Range("C48").Value = 'PKI'
$ws.Range("C49").Value = 'A web server facing the Internet, an application server on the internal network, a database server on the internal network'
$ws.Range("C50").Value = 'RSA is asymmetric, which is used to create a public/private key pair; AES is symmetric, which is used to encrypt data'
$ws.Range("C51").Value = 'tcptrace'
$ws.Range("C52").Value = 'T5 '
$ws.Range("C53").Value = 'tcpdump'
$ws.Range("C55").Value = 'Social engineering'
$ws.Range("C56").Value = 'NTP'
$ws.Range("C57").Value = 'tcp.port ==21'
$ws.Range("C58").Value = 'Encrypt the backup tapes and transport them in a lockbox'
$ws.Range("C59").Value = 'Results about all discounts from the site eccouncil.org except for the ilearn format'
$ws.Range("C60").Value = 'The tester only partially knows the internal structure'
$ws.Range("C61").Value = 'The data collection speeds, data processing speed, or enrolment time'
$ws.Range("C62").Value = 'He is scanning from 192.168.1.64 to 192.168.1.78 because of the mask /28 and the servers are not in that range'
$ws.Range("C63").Value = 'An inability to access any website'
$ws.Range("C64").Value = 'Social Engineering'
$ws.Range("C65").Value = 'Chosen-plaintext'

# Switch active sheet/tab to sheet "3" and set final selection
$ws.Activate()
$ws.Range("C66").Select()
